$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.967.69'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '2.540.45'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.08'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.34'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.64'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.112'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '2.933.74'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('D15').Value = '2.543.16'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.41'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.851'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '43.051.57'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.09'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.65'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.30%  '
$ws.Range('D21').Value = '0.0₃0972'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.45'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '252.72'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.43'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.08'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.27'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.07'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.33'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.08'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0792'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.112'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.92%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.97'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('E41').Value = '  +10.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.82'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0305'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.34'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = '2.026.62'
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.67'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.82'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.90%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.787.04'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.84'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.56'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.71%  '
